$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to match repulled data
$ws.Range("F11").Value = 2
$ws.Range("F17").Value = -8
$ws.Range("F20").Value = -1
$ws.Range("F27").Value = -4
$ws.Range("F28").Value = -1
$ws.Range("F37").Value = -2
$ws.Range("F39").Value = -4
$ws.Range("F42").Value = -1
$ws.Range("F46").Value = 1
$ws.Range("F47").Value = 5
$ws.Range("F49").Value = 4
$ws.Range("F50").Value = 3
$ws.Range("F51").Value = 2
$ws.Range("F53").Value = 0
$ws.Range("F60").Value = -2
$ws.Range("F66").Value = 0
$ws.Range("F69").Value = 4
$ws.Range("F72").Value = -6
